$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Slpi"
$ws.Range("C2").Value = "Plscr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7680743333333333
$ws.Range("H2").Value = 2.304223
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.581756333333333
$ws.Range("N2").Value = 7.745269
$ws.Range("O2").Value = 0.3831594293522085
$ws.Range("P2").Value = 0.3831594293522085
$ws.Range("Q2").Value = 1.982980774554111
$ws.Range("R2").Value = 17.846826970987
$ws.Range("S2").Value = 0.3831594293522085
$ws.Range("T2").Value = 0.3831594293522085

$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Slpi"
$ws.Range("C3").Value = "Plscr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7680743333333333
$ws.Range("H3").Value = 2.304223
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.786435
$ws.Range("N3").Value = 11.359305
$ws.Range("O3").Value = 0.561946243782842
$ws.Range("P3").Value = 0.5619462437828419
$ws.Range("Q3").Value = 2.908263538335
$ws.Range("R3").Value = 26.174371845015
$ws.Range("S3").Value = 0.561946243782842
$ws.Range("T3").Value = 0.5619462437828419

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Slpi"
$ws.Range("C4").Value = "Plscr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7680743333333333
$ws.Range("H4").Value = 2.304223
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.369882
$ws.Range("N4").Value = 1.109646
$ws.Range("O4").Value = 0.05489432686494953
$ws.Range("P4").Value = 0.05489432686494952
$ws.Range("Q4").Value = 0.284096870562
$ws.Range("R4").Value = 2.556871835058
$ws.Range("S4").Value = 0.05489432686494953
$ws.Range("T4").Value = 0.05489432686494952

$ws.Rows("5:7").Delete() | Out-Null
